$d = $word.ActiveDocument
$dash = [char]0x2013

# ---------------------------------------------------------------
# Change 1: "There are four FPGAs on the FEB..." paragraph gets a
# part-number parenthetical inserted after "FEB".
# ---------------------------------------------------------------
$r = $d.Content
$found = $r.Find.Execute("FPGAs on the FEB.", $true, $false, $false, $false, $false,
                          $true, 1, $false, "", 0)
if (-not $found) { throw "Could not find 'FPGAs on the FEB.' anchor" }

# Collapse to the point right before the period that follows "FEB".
$insPos = $r.End - 1
$insPoint = $d.Range($insPos, $insPos)
$insPoint.InsertBefore(" (part number xc7s50fgga484-2)")

# " (" -- plain formatting (already inherited, nothing to change)
# "part number" -- Calibri 11pt
$pnStart = $insPos + 2
$pnEnd = $pnStart + "part number".Length
$pnRange = $d.Range($pnStart, $pnEnd)
$pnRange.Font.Name = "Calibri"
$pnRange.Font.Size = 11

# " xc7s50fgga484-2" -- Calibri 11pt
$partStart = $pnEnd
$partEnd = $partStart + " xc7s50fgga484-2".Length
$partRange = $d.Range($partStart, $partEnd)
$partRange.Font.Name = "Calibri"
$partRange.Font.Size = 11

# ")" -- Calibri 11pt
$closeStart = $partEnd
$closeEnd = $closeStart + ")".Length
$closeRange = $d.Range($closeStart, $closeEnd)
$closeRange.Font.Name = "Calibri"
$closeRange.Font.Size = 11

# ---------------------------------------------------------------
# Changes 2-5: the four "0x### .. 0x### - FPGAn. CMB Inputs 1..4 "
# lines get their split runs (and proofErr markers) merged back
# into a single run. The visible text is unchanged, so a
# find/replace of the text with itself collapses the runs.
# ---------------------------------------------------------------
$ranges = @(
    "0x000..0x2FF $dash FPGA1. CMB Inputs 1..4 ",
    "0x400..0x6FF $dash FPGA2. CMB Inputs 1..4 ",
    "0x800..0xAFF $dash FPGA3. CMB Inputs 1..4 ",
    "0xC00..0xEFF $dash FPGA4. CMB Inputs 1..4 "
)
foreach ($t in $ranges) {
    $r2 = $d.Content
    $ok = $r2.Find.Execute($t, $true, $false, $false, $false, $false,
                            $true, 1, $false, $t, 2)
    if (-not $ok) { throw "Could not find/replace '$t'" }
}

# ---------------------------------------------------------------
# Change 6a: 0x022 paragraph - merge "upper 16 bits of the " +
# "32 bit" + " test counter. " into a single run (proofErr removed).
# ---------------------------------------------------------------
$t1 = "upper 16 bits of the 32 bit test counter. "
$r3 = $d.Content
$ok = $r3.Find.Execute($t1, $true, $false, $false, $false, $false,
                        $true, 1, $false, $t1, 2)
if (-not $ok) { throw "Could not find/replace upper-bits sentence" }

# ---------------------------------------------------------------
# Change 6b: 0x023 paragraph - merge "lower 16 bits of the " +
# "32 bit" + " test counter.  " into a single run, and change the
# following sentence from "A read from this address displays..."
# to "A read to this address displays...".
# ---------------------------------------------------------------
$t2 = "lower 16 bits of the 32 bit test counter.  "
$r4 = $d.Content
$ok = $r4.Find.Execute($t2, $true, $false, $false, $false, $false,
                        $true, 1, $false, $t2, 2)
if (-not $ok) { throw "Could not find/replace lower-bits sentence" }

$r5 = $d.Content
$ok = $r5.Find.Execute("A read from this address displays the value of the ",
                        $true, $false, $false, $false, $false,
                        $true, 1, $false,
                        "A read to this address displays the value of the ", 2)
if (-not $ok) { throw "Could not find/replace 'A read from this address...'" }

Write-Output "done"
